$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, pushing existing rows 2,3,4 down to 3,4,5
$ws.Rows.Item(2).Insert()

# Fill the new row 2 with kebab-case slug versions of the header row labels
$ws.Range("A2").Value = "municipio-codigo"
$ws.Range("B2").Value = "estado-civil"
$ws.Range("C2").Value = "estado-civil-codigo"
$ws.Range("D2").Value = "sexo-codigo"
$ws.Range("E2").Value = "sexo"
$ws.Range("F2").Value = "municipio-nombre"
